# Apply new built-in table style to the three tables on slides 14-16
$p = $ppt.ActivePresentation
$newTableStyle = "{D2159ADB-F392-4E4C-9ED8-C435E1DE8964}"
foreach ($slideIndex in 14,15,16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    $table = $tableShape.Table
    $table.ApplyStyle($newTableStyle)
}

# Switch the deck's master theme palette from the custom 'Integral / Red Violet'
# colours over to the stock Office theme palette (Design > Colors)
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
$themeColors.Colors(1).RGB = 0 # dk1 -> #000000
$themeColors.Colors(2).RGB = 16777215 # lt1 -> #FFFFFF
$themeColors.Colors(3).RGB = 6968388 # dk2 -> #44546A
$themeColors.Colors(4).RGB = 15132391 # lt2 -> #E7E6E6
$themeColors.Colors(5).RGB = 13998939 # accent1 -> #5B9BD5
$themeColors.Colors(6).RGB = 3243501 # accent2 -> #ED7D31
$themeColors.Colors(7).RGB = 10855845 # accent3 -> #A5A5A5
$themeColors.Colors(8).RGB = 49407 # accent4 -> #FFC000
$themeColors.Colors(9).RGB = 12874308 # accent5 -> #4472C4
$themeColors.Colors(10).RGB = 4697456 # accent6 -> #70AD47
$themeColors.Colors(11).RGB = 12673797 # hlink -> #0563C1
$themeColors.Colors(12).RGB = 7491477 # folHlink -> #954F72
